$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 8 & 9 values (order chosen to reproduce the shared-string table order) ---
$ws.Range("A8").Value = "2_3"
$ws.Range("B8").Value = "Create new Textview and spinner dynamically"
$ws.Range("D8").Value = "Yes"
$ws.Range("D9").Value = "No"
$ws.Range("C8").Value = "When user click button on Add New Word, EditText and spinner are added to`nmore item add"
$ws.Range("E8").Value = "Create new Linear Layout and `nput new dynamically created Spinner and EditText on the layout."
$ws.Range("F8").Value = "AddNewWordDialog`ndialog_add_new_word.xml"
$ws.Range("A9").Value = "2_4"
$ws.Range("B9").Value = "Create Buttons on Add word dialog"
$ws.Range("C9").Value = "Create Add button to save the word into database and cancel button to cancel the work"

Write-Output "done"
